$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update existing "Name of file" / filename -> "File" / generic csv filename
$ws.Range("A2").Value = "File"
$ws.Range("B2").Value = "scopus_social_evol_2018.csv"
$ws.Range("C2").Value = "filename"

# New rows 3-6 describing the generic file config
$ws.Range("A3").Value = "ID column"
$ws.Range("B3").Value = "ID-1"

$ws.Range("A4").Value = "Title"
$ws.Range("B4").Value = "Title of the contribution in original language-2"

$ws.Range("A5").Value = "ISSN"
$ws.Range("B5").Value = "Journal > ISSN-5"

$ws.Range("A6").Value = "DOI"
$ws.Range("B6").Value = "Electronic version(s) of this work > DOI (Digital Object Identifier)-6"

# Selection moves to B3 as in the diff
$ws.Range("B3").Select()
